# "Generate Report for Handback"
#
# Fills in the handback status/columns on the Overview, zh-cn and de-de
# sheets now that the localized content has come back in sync with en-US.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$mdName   = "d3627a0d-745b-4ce1-9ae3-92e25b28d728.md"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0fb90ebf682bd75a07347823630fae9f51a67a17/e2e/d3627a0d-745b-4ce1-9ae3-92e25b28d728.md"

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns (zh-cn / de-de) now show "handed back" ---
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("I2").Value = $mdName
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdName)
$zhcn.Range("J2").Value = "d3627a0d-745b-4ce1-9ae3-92e25b28d728.638d068fd090a30f31bc3a8c579211c0518d986c.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 17:09:44"

# --- de-de sheet ---
$dede.Range("C2").Value = $handedBack
$dede.Range("I2").Value = $mdName
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdName)
$dede.Range("J2").Value = "d3627a0d-745b-4ce1-9ae3-92e25b28d728.638d068fd090a30f31bc3a8c579211c0518d986c.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 17:09:52"
